$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing data rows (rows 3-6), shifting nothing below them up.
$ws.Range("A3:E6").EntireRow.Delete()

# Update the remaining data row's process/machine columns.
$ws.Cells.Item(2, 4).Value = "قطع ليزر"
$ws.Cells.Item(2, 5).Value = "ماكينة قطع باليزر 3 م"
